$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new "sex" column (D) with "U" for rows 2-5
$ws.Range("D2:D5").Value = "U"

# Fill in the new "environ" column (E) with 1 for rows 2-5
$ws.Range("E2:E5").Value = 1

# Fill in the "sire" column (F) for rows 4-5
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(5, 6).Value = 3

# Fill in the "dam" column (G) for rows 4-5
$ws.Cells.Item(4, 7).Value = 2
$ws.Cells.Item(5, 7).Value = 3

# Update the selection to D2:D5 with active cell D2
$ws.Range("D2:D5").Select()
